# Hoverbot alert sound implemented
#
# The HoverbotAlert row (row 2, "Sound when the hoverbot notices the
# player") is marked Completed and a follow-up note is added asking for
# the sound to be made more hostile-sounding.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Completed"
$ws.Range("F2").Value = "Might need to be made a bit more hostile sounding"

# Leave the selection on the freshly-entered note, matching where the
# author's cursor was when they last saved the file.
$ws.Range("F2").Select()
